$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 133.4
$ws.Range("B3").Value = 178.6
$ws.Range("C3").Value = 161.9
$ws.Range("C6").Value = 64.5
